$d = $word.ActiveDocument

function Find-ParagraphLike($pattern) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $pattern) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. Skills list: remove the standalone "Android Studio Development"
#    bullet and rewrite the "Fluent in: ..." bullet that follows it.
# ---------------------------------------------------------------------
$pAndroid = Find-ParagraphLike "Android Studio Development*"
if ($pAndroid -ne $null) {
    $pAndroid.Range.Delete()
}

$pFluent = Find-ParagraphLike "Fluent in:*"
if ($pFluent -ne $null) {
    $r = $pFluent.Range
    $r.Find.Execute("Fluent in: Python, PowerShell, C++, Java, and Golang", $true, $false, $false, $false, $false, $true, 1, $false, "Fluent in: Golang, Shell-Scripting, HTML-CSS, and Javascript", 2)
}

# ---------------------------------------------------------------------
# 2. Cloud infra bullet: "...mainly AWS)." -> "...mainly AWS and EKS
#    Kubernetes)."
# ---------------------------------------------------------------------
$pAws = Find-ParagraphLike "Deploying/maintaining cloud infrastructures*"
if ($pAws -ne $null) {
    $r = $pAws.Range
    $r.Find.Execute("mainly AWS).", $true, $false, $false, $false, $false, $true, 1, $false, "mainly AWS and EKS Kubernetes).", 2)
}

# ---------------------------------------------------------------------
# 3. "Web Developer" job -> "Cloud Developer" job at AT&T, 2021, ending
#    2/02/2021; body paragraph rewritten.
# ---------------------------------------------------------------------
$pWebHeading = Find-ParagraphLike "Web Developer*Self*2017*Present*"
if ($pWebHeading -ne $null) {
    $r = $pWebHeading.Range
    $r.Find.Execute("Web Developer", $true, $false, $false, $false, $false, $true, 1, $false, "Cloud Developer", 2)
    $r = $pWebHeading.Range
    $r.Find.Execute("Self", $true, $false, $false, $false, $false, $true, 1, $false, "AT&T", 2)
    $r = $pWebHeading.Range
    $r.Find.Execute("2017", $true, $false, $false, $false, $false, $true, 1, $false, "2021", 2)
    $r = $pWebHeading.Range
    $r.Find.Execute("Present", $true, $false, $false, $false, $false, $true, 1, $false, "2/02/2021", 2)
}

$pWebBody = Find-ParagraphLike "Designed and launched websites*"
if ($pWebBody -ne $null) {
    $r = $pWebBody.Range
    $old = "Designed and launched websites/projects with full-stack development patterns. Utilized the following for backend: Google Dev Console, AWS, Python, Java, and Golang. For Frontend: Javascript, HTML, CSS, and Vue."
    $new = "Worked with  CI-CD/SRE teams to help develop code pipelines and develop a Network Cloud for AT&T. Delivered test-suites and pipelines utilizing Jenkins/Groovy deployed to Kubernetes cloud environements."
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# ---------------------------------------------------------------------
# 4. "SuperDBWebApp" project -> "The-Informer(Alpha-development)",
#    date 4/5/20 -> 11/15/20; body paragraph rewritten.
# ---------------------------------------------------------------------
$pSuperHeading = Find-ParagraphLike "SuperDBWebApp*"
if ($pSuperHeading -ne $null) {
    $r = $pSuperHeading.Range
    $r.Find.Execute("SuperDBWebApp", $true, $false, $false, $false, $false, $true, 1, $false, "The-Informer(Alpha-development)", 2)
    $r = $pSuperHeading.Range
    $r.Find.Execute("4/5/20", $true, $false, $false, $false, $false, $true, 1, $false, "11/15/20", 2)
}

$pSuperBody = Find-ParagraphLike "A website/dat*"
if ($pSuperBody -ne $null) {
    $r = $pSuperBody.Range
    $old = "A website/database hybrid running on a Docker container, utilizing Amazon MongoDB services to host files/data. (Currently in testing, see my website for link.)"
    $new = "A project made with Kubernetes deployed Microservices in AWS. Utilizes text message APIs as a platform for charity organizations and campaigns to reach constituents and combat misinformation online."
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done"
